$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 19:40"

# --- Reorder Irlanda to right after Libia (row 69), shifting Ghana and
#     Estado de Palestina down one row, and refresh Irlanda's stats.
#     Row 69 = Libia (unchanged). Rows 70-72 get rewritten in place:
#       new row 70 = Irlanda (updated numbers)
#       new row 71 = Ghana   (same numbers it had before the move)
#       new row 72 = Estado de Palestina (same numbers it had before the move)

$ws.Range("A70").Value = "Irlanda"
$ws.Range("B70").Value = 47427
$ws.Range("C70").Value = 998
$ws.Range("D70").Value = 23364
$ws.Range("E70").Value = 22222
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 1841

$ws.Range("A71").Value = "Ghana"
$ws.Range("B71").Value = 47173
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 46527
$ws.Range("E71").Value = 336
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 310

$ws.Range("A72").Value = "Estado de Palestina"
$ws.Range("B72").Value = 46434
$ws.Range("C72").Value = 334
$ws.Range("D72").Value = 39921
$ws.Range("E72").Value = 6111
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 402

# --- General covid-numbers refresh for other countries ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 8235136
$ws.Range("C4").Value = 18821
$ws.Range("E4").Value = 2682921
$ws.Range("G4").Value = 328
$ws.Range("H4").Value = 223045

# Row 5: India
$ws.Range("B5").Value = 7421141
$ws.Range("C5").Value = 55632
$ws.Range("D5").Value = 6510749
$ws.Range("E5").Value = 797472
$ws.Range("G5").Value = 774
$ws.Range("H5").Value = 112920

# Row 13: Francia
$ws.Range("B13").Value = 834770
$ws.Range("C13").Value = 25086
$ws.Range("E13").Value = 697441
$ws.Range("G13").Value = 122
$ws.Range("H13").Value = 33247

# Row 21: Alemania
$ws.Range("B21").Value = 355010
$ws.Range("C21").Value = 6194
$ws.Range("E21").Value = 60578
$ws.Range("G21").Value = 22
$ws.Range("H21").Value = 9832

# Row 24: Turquia
$ws.Range("B24").Value = 343955
$ws.Range("C24").Value = 1812
$ws.Range("D24").Value = 301098
$ws.Range("E24").Value = 33704
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 9153

# Row 64: Libano
$ws.Range("B64").Value = 60113
$ws.Range("C64").Value = 1368
$ws.Range("D64").Value = 26468
$ws.Range("E64").Value = 33136
$ws.Range("G64").Value = 8
$ws.Range("H64").Value = 509

# Row 66: Argelia
$ws.Range("B66").Value = 53998
$ws.Range("C66").Value = 221
$ws.Range("D66").Value = 37856
$ws.Range("E66").Value = 14301
$ws.Range("G66").Value = 14
$ws.Range("H66").Value = 1841

# Row 90: Republica de Macedonia
$ws.Range("B90").Value = 22607
$ws.Range("C90").Value = 437
$ws.Range("D90").Value = 16949
$ws.Range("E90").Value = 4837
$ws.Range("G90").Value = 6
$ws.Range("H90").Value = 821

# Row 165: Liberia
$ws.Range("B165").Value = 1377
$ws.Range("C165").Value = 3
$ws.Range("D165").Value = 1264
$ws.Range("E165").Value = 31

# Row 189: Monaco
$ws.Range("B189").Value = 255
$ws.Range("C189").Value = 2
$ws.Range("E189").Value = 36
